# Weekly price-sheet update: a new record is published for this market/product
# combination each week. The new record is inserted as row 155 (pushing every
# existing record down by one row), so the oldest record that used to live in
# row 238 now overflows into a brand-new row 239.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 155, shifting rows 155:238 down to 156:239.
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A155").Value = 7
$ws.Range("B155").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C155").Value = "Ñuble"
$ws.Range("D155").Value = 44572
$ws.Range("E155").Value = 16
$ws.Range("F155").Value = 100114013
$ws.Range("G155").Value = "Zanahoria"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 120
$ws.Range("K155").Value = 6000
$ws.Range("L155").Value = 6500
$ws.Range("M155").Value = 6250
$ws.Range("N155").Value = "`$/saco 20 kilos"
$ws.Range("O155").Value = "Provincia de Diguillín"
$ws.Range("P155").Value = 312
$ws.Range("Q155").Value = 20
$ws.Range("R155").Value = "Hortaliza"
